$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.130.86'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.37%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.314.57'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.97%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.97%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.47'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.67%  '

$ws.Range("E7").Value = '  +1.85%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("E9").Value = '  +5.11%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.14'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +9.28%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0793'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.86%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.117'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.14%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '17.98'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +13.73%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.88'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.26%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.676.11'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.00%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.346.02'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.05%  '

$ws.Range("E17").Value = '  +3.16%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.012.79'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.12%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.63'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +8.06%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.22'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.13%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0903'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.36%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.84'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.44%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.18'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.53%  '

$ws.Range("E24").Value = '  +12.84%  '

$ws.Range("E25").Value = '  +0.46%  '

$ws.Range("E27").Value = '  +4.26%  '

$ws.Range("E28").Value = '  +3.14%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34.72'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.94%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '168.60'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.66%  '

$ws.Range("E31").Value = '  +0.56%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.01%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.02'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.70%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.72'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.24%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.47'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.08%  '

$ws.Range("E36").Value = '  +2.23%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0693'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.71%  '

$ws.Range("E38").Value = '  +4.37%  '

$ws.Range("E39").Value = '  +2.10%  '

$ws.Range("E40").Value = '  +4.05%  '

$ws.Range("E41").Value = '  +1.26%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.985.21'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.37%  '

$ws.Range("E43").Value = '  +4.67%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.25'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.70%  '

$ws.Range("E45").Value = '  +7.34%  '

$ws.Range("E46").Value = '  +5.20%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '17.55'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.04%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '56.17'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +8.09%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.543.70'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.90%  '

$ws.Range("E50").Value = '  +3.96%  '

$ws.Range("B51").Value = 'HuobiToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.75'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.91%  '
